# Re-write the Login page / user list:
#   - Insert a "group" (分组) column and an "admin account" (管理员账号)
#     column between the existing "username" and "random password" columns.
#   - Add two new users (谢江霞, 孙洪莹) with their generated passwords.
#   - Fill in the group/role for every user, and mark the admin account.
#
# The existing B (random password) / L (label) / M (generator formula)
# columns are pushed two columns to the right by the column insert below,
# which is exactly what happens when two new columns are inserted in the
# Excel UI - the original cells (and their shared-string references) keep
# their contents and simply slide from B->D and L/M->N/O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new, empty columns at B/C - existing data in B (and
#    everything to its right) shifts right by two columns.
$ws.Range("B1:C1").EntireColumn.Insert() | Out-Null

# 2. New rows for the two newly added accounts. Username + password are
#    filled in first so the new shared strings line up the same way the
#    source workbook orders them.
$ws.Range("A8").Value = "谢江霞"
$ws.Range("A9").Value = "孙洪莹"
$ws.Range("D8").Value = "O8q93598"
$ws.Range("D9").Value = "N5o49899"

# 3. New header cells for the inserted columns.
$ws.Range("B1").Value = "分组"
$ws.Range("C1").Value = "管理员账号"

# 4. Flag the one admin account.
$ws.Range("C4").Value = "是"

# 5. Fill in every user's group / role.
$ws.Range("B2").Value = "地区经理"
$ws.Range("B3").Value = "地区经理"
$ws.Range("B4").Value = "全国总监"
$ws.Range("B5").Value = "大区总监"
$ws.Range("B6").Value = "大区经理"
$ws.Range("B7").Value = "大区总监"
$ws.Range("B8").Value = "地区经理"
$ws.Range("B9").Value = "地区经理"

# 6. Leave the selection where the edit naturally finished.
$ws.Range("F6").Select() | Out-Null
